# Updates cryptos price/volume columns (D, E) for rows 2-51 per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.100.95"
$ws.Range("E2").Value = "'  +0.21%  "
$ws.Range("D3").Value = "'3.834.11"
$ws.Range("E3").Value = "'  -0.27%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'706.59"
$ws.Range("E5").Value = "'  +2.02%  "
$ws.Range("D6").Value = "'171.20"
$ws.Range("E6").Value = "'  -0.92%  "
$ws.Range("D7").Value = "'3.831.42"
$ws.Range("E7").Value = "'  -0.20%  "
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("E9").Value = "'  -0.42%  "
$ws.Range("E10").Value = "'  -0.92%  "
$ws.Range("E11").Value = "'  +0.64%  "
$ws.Range("E12").Value = "'  -0.52%  "
$ws.Range("E13").Value = "'  -1.46%  "
$ws.Range("D14").Value = "'36.53"
$ws.Range("E14").Value = "'  -0.21%  "
$ws.Range("D15").Value = "'4.481.19"
$ws.Range("E15").Value = "'  -0.27%  "
$ws.Range("D16").Value = "'3.773.30"
$ws.Range("E16").Value = "'  -2.57%  "
$ws.Range("D17").Value = "'71.179.19"
$ws.Range("E17").Value = "'  +0.19%  "
$ws.Range("E18").Value = "'  -0.34%  "
$ws.Range("E19").Value = "'  +0.16%  "
$ws.Range("D20").Value = "'17.37"
$ws.Range("E20").Value = "'  -2.51%  "
$ws.Range("D21").Value = "'495.31"
$ws.Range("E21").Value = "'  +1.53%  "
$ws.Range("E22").Value = "'  -4.96%  "
$ws.Range("E23").Value = "'  +2.09%  "
$ws.Range("D24").Value = "'85.64"
$ws.Range("E24").Value = "'  +1.12%  "
$ws.Range("E25").Value = "'  -1.64%  "
$ws.Range("D26").Value = "'10.61"
$ws.Range("E26").Value = "'  +1.09%  "
$ws.Range("D27").Value = "'12.11"
$ws.Range("E27").Value = "'  -2.22%  "
$ws.Range("D28").Value = "'3.987.71"
$ws.Range("E28").Value = "'  -0.40%  "
$ws.Range("E29").Value = "'  -2.95%  "
$ws.Range("E30").Value = "'  +0.00%  "
$ws.Range("E31").Value = "'  -0.82%  "
$ws.Range("E32").Value = "'  -2.70%  "
$ws.Range("D33").Value = "'2.23"
$ws.Range("E33").Value = "'  -2.77%  "
$ws.Range("E34").Value = "'  -1.26%  "
$ws.Range("E35").Value = "'  -2.80%  "
$ws.Range("D36").Value = "'3.802.77"
$ws.Range("E36").Value = "'  +0.14%  "
$ws.Range("D37").Value = "'9.15"
$ws.Range("E37").Value = "'  -1.12%  "
$ws.Range("E38").Value = "'  -0.31%  "
$ws.Range("E39").Value = "'  -1.44%  "
$ws.Range("E40").Value = "'  +4.00%  "
$ws.Range("E41").Value = "'  -2.70%  "
$ws.Range("E42").Value = "'  -1.06%  "
$ws.Range("D43").Value = "'3.33"
$ws.Range("E43").Value = "'  -3.24%  "
$ws.Range("E45").Value = "'  -0.16%  "
$ws.Range("E46").Value = "'  +1.40%  "
$ws.Range("D47").Value = "'163.98"
$ws.Range("E47").Value = "'  -0.33%  "
$ws.Range("D48").Value = "'429.86"
$ws.Range("E48").Value = "'  +3.99%  "
$ws.Range("D49").Value = "'48.94"
$ws.Range("E49").Value = "'  +0.50%  "
$ws.Range("D50").Value = "'8.75"
$ws.Range("E50").Value = "'  +0.74%  "
$ws.Range("E51").Value = "'  -1.74%  "
